# updated LLY, ABBV, AMD
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Bump "last updated" dates for the affected screens:
#  D5  -> Biopharma screen (covers LLY, ABBV, etc.)   03/14/23 -> 05/02/23
#  D10 -> Hardware & Semis screen (covers AMD, etc.)  02/04/23 -> 05/02/23
#  D11 -> Software screen                             04/25/23 -> 05/01/23
$ws.Range("D5").Value = 45048
$ws.Range("D10").Value = 45048
$ws.Range("D11").Value = 45047

# Move the active selection to D6, matching where the editor left off.
$ws.Range("D6").Select()
